$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

# Row 17: new customer "lucas" appended below the last existing row (16)
$ws.Range("A17").Value = "lucas"
$ws.Range("B17").Value = "asdasdas"
$ws.Range("C17").Value = "asdasdas"
# D17 intentionally left blank (matches source data - no address column D value)
$ws.Range("E17").Value = "92320-195"
$ws.Range("F17").Value = "joanues@gmail.com"

# G17 holds a long, purely numeric phone string; format the cell as Text first
# so Excel keeps it as a text value instead of converting it into a number
# (mirrors every other row in this column, which is stored as text).
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "12312312312"

$ws.Range("H17").Value = "Rua 3 Pinheiros I, 27"
